$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells stay text (matching original inlineStr typing) by setting
# the number format to Text before assigning numeric-looking strings, so
# Excel does not silently convert values like "7.70" -> 7.7.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.771.84"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.365.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.70"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.70"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.14%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.939.26"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.86"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.356.69"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.899.86"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.47"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.89"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.06"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.07"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.03%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.189"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.01%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.10"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.86"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -7.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.94"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "167.04"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.400.25"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.44"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.54"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.81%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.89%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.441.61"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.57%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.61"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.08%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -6.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0256"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.95"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -6.16%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.87%  "
